# testng dataprovider examples added
#
# Rewrites the small username/password example sheet into a set of
# TestNG @DataProvider style example rows, and highlights the header
# row ("username" / "pass") using Excel's built-in "Good" cell style.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row -----------------------------------------------------
$ws.Range("A1").Value = "username"
$ws.Range("B1").Value = "pass"
$ws.Range("A1:B1").Style = "Good"

# --- Data rows --------------------------------------------------------
$ws.Range("A2").Value = "one"
$ws.Range("B2").Value = "a"

$ws.Range("A3").Value = "two"
$ws.Range("B3").Value = "abcde"

$ws.Range("A4").Value = "three"
$ws.Range("B4").Value = "a"

# --- Column sizing: widen column A so the longer sample values fit ---
$ws.Columns("A").ColumnWidth = 16.33

# --- Leave the selection where the author last left it ---------------
$ws.Range("H9:H10").Select() | Out-Null
